$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" - the localization-status report is updated
# once a handback (target -> source sync) completes:
#   * Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#   * The per-language "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated with the handback
#     artifacts that were just produced.
# ---------------------------------------------------------------------------

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Update the status text everywhere it appears -------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Hyperlink target addresses (same repo blob used by column A) ---------
$addrMain = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61b9e4dbe755eaf586d8d5858b933dfabeced4d6/e2e/43ba885b-911d-48b2-847f-0c9ff3369a15.md"
$addrFfff = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61b9e4dbe755eaf586d8d5858b933dfabeced4d6/e2e/ffff6cfd71fe-2663-4174-9774-60d3076beb11.md"
$mainName = "43ba885b-911d-48b2-847f-0c9ff3369a15.md"
$ffffName = "ffff6cfd71fe-2663-4174-9774-60d3076beb11.md"

# --- zh-cn: populate "Latest Target File" (I), "Latest Handback File" (J),
#     and "Latest Handback DateTime" (K) for both data rows -----------------
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("J2").Value = "43ba885b-911d-48b2-847f-0c9ff3369a15.1788a6d4d7b073fbd3126a615837aab506dbce54.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-03 17:10:49"
$wsZhCn.Range("J3").Value = "43ba885b-911d-48b2-847f-0c9ff3369a15.1788a6d4d7b073fbd3126a615837aab506dbce54.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-03 17:10:49"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $addrFfff, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $ffffName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)

# --- de-de: same shape, but its own handback xlf + a later timestamp ------
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("J2").Value = "43ba885b-911d-48b2-847f-0c9ff3369a15.1788a6d4d7b073fbd3126a615837aab506dbce54.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 17:10:57"
$wsDeDe.Range("J3").Value = "43ba885b-911d-48b2-847f-0c9ff3369a15.1788a6d4d7b073fbd3126a615837aab506dbce54.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 17:10:57"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $addrFfff, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $ffffName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $addrMain, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mainName)

# --- Widen the columns that now hold the longer status / file-name text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40
